$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "286.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.55%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.38"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.77%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.108"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.64%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06662"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.58%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.332"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.31%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.403"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.32%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.346"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.69%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9205"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.40%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1563"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.61%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06476"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.29%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07570"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.49%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02905"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.20%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.08995"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.15%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001592"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.49%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04479"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.45%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0006463"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.82%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006298"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.30%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.456"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.07%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.24%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.74%"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.84%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.062"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.22%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1549"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.96%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001190"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.39%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004121"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-5.07%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "5.99%"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.04%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04181"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.58%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006734"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.27%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1239"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-12.46%"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.36%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01264"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "6.92%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005620"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.45%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.01307"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-29.30%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.967"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "26.04%"
